# Planning Sheet update: "home page of flight booking system"
#
# 1. Remove the "Actual ETA" column (column H) entirely - STATUS and the
#    remarks column shift left to take its place.
# 2. Renumber the "Task no." column sequentially (1-12) by row position.
# 3. Rename the last header from "STATEMENTS" to "Remarks/Comments".
# 4. Fill in / fix a few remark cells (typo fix, new remark, cleared remark).
# 5. Push a couple of task dates out by a day (rows 10 and 12).
# 6. Move the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the "Actual ETA" column (H). Cells in I/J shift to H/I. ---
$ws.Columns.Item(8).Delete()

# --- 2. Renumber Task no. column (A) for rows 2-13. ---
$taskNumbers = @{
    2 = 1
    3 = 2
    4 = 3
    5 = 4
    6 = 5
    7 = 6
    8 = 7
    9 = 8
    10 = 9
    11 = 10
    12 = 11
    13 = 12
}
foreach ($row in $taskNumbers.Keys) {
    $ws.Cells.Item($row, 1).Value = $taskNumbers[$row]
}

# --- 3. Header rename: column I (after the shift) is now "Remarks/Comments". ---
$ws.Cells.Item(1, 9).Value = "Remarks/Comments"

# --- 4. Remark fixups in column I. ---
$ws.Cells.Item(2, 9).Value = "Login to Home connection database connectivity"
$ws.Cells.Item(4, 9).Value = "Merging and connected to home page form left"
$ws.Cells.Item(7, 9).Value = ""

# --- 5. Date shifts. ---
$ws.Cells.Item(10, 6).Value = "2023-05-26"
$ws.Cells.Item(10, 7).Value = "2023-05-27"
$ws.Cells.Item(12, 6).Value = "2023-05-27"

# --- 6. Selection moves to H19. ---
$ws.Range("H19").Select()
